# Apply the edits described in the diff: four paragraphs that used to
# contain run-on sentences (several ideas glued together with no
# separator) are split so each idea sits on its own line, separated by
# a manual line break (<w:br/>) inside the same run. We reproduce this
# with Find/Replace using the "^l" special code, which Word's find
# engine turns into a manual line break when used as replacement text.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, `
                                      $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Could not find text: $old"
    }
}

# 1) Objetivos paragraph -> split into three lines
Replace-Text `
    "- Capacitar o aluno para relacionar as propriedades químicas e físicas dos elementos e seus compostos com suas posições na tabela periódica.-Capacitar o aluno a escrever os métodos industriais de obtenção dos elementos e seus compostos, bem como descrever suas aplicações- capacitar o aluno a comunicar-se eficazmente nas formas escrita, oral e gráfica" `
    "- Capacitar o aluno para relacionar as propriedades químicas e físicas dos elementos e seus compostos com suas posições na tabela periódica.^l-Capacitar o aluno a escrever os métodos industriais de obtenção dos elementos e seus compostos, bem como descrever suas aplicações^l- capacitar o aluno a comunicar-se eficazmente nas formas escrita, oral e gráfica"

# 2) Programa resumido paragraph -> split into two lines
Replace-Text `
    "- Metais Representativos: Metais do Grupo 1, Metais do Grupo 2 e Metais do Grupo 13.- Metais de Transição: Propriedades gerais, Complexos." `
    "- Metais Representativos: Metais do Grupo 1, Metais do Grupo 2 e Metais do Grupo 13.^l- Metais de Transição: Propriedades gerais, Complexos."

# 3) Programa paragraph -> split into two lines
Replace-Text `
    "Metais e compostos dos grupos 1, 2, 13 e de transição da Tabela Periódica: Propriedades físicas e químicas (relação com a posição na Tabela Periódica), processos de obtenção dos metais e compostos e aplicações - Formação de Complexos.Relacionar a disciplina com disciplinas anteriores e posteriores da grade do curso." `
    "Metais e compostos dos grupos 1, 2, 13 e de transição da Tabela Periódica: Propriedades físicas e químicas (relação com a posição na Tabela Periódica), processos de obtenção dos metais e compostos e aplicações - Formação de Complexos.^lRelacionar a disciplina com disciplinas anteriores e posteriores da grade do curso."

# 4) Avaliação / Método run -> split into two lines
Replace-Text `
    "A avaliação tem como requisito quantificar as competências adquiridas conforme objetivadas.Duas provas escritas (P1 e P2) e listas de exercícios de acompanhamento continuado. A partir das notas das listas de exercício será calculada a média, LE." `
    "A avaliação tem como requisito quantificar as competências adquiridas conforme objetivadas.^lDuas provas escritas (P1 e P2) e listas de exercícios de acompanhamento continuado. A partir das notas das listas de exercício será calculada a média, LE."
